$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the existing _GoBack bookmark from the end of the "average
#    campaign times" paragraph. It will be relocated (see step 2) to mirror
#    the author's last edit position after making the below corrections.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Split the first (bold) heading paragraph right after "can " and drop a
#    fresh _GoBack bookmark at that split point.
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute("make about Kickstarter campaigns given the provided data?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $d.Range($rng1.Start, $rng1.Start)
$d.Bookmarks.Add("_GoBack", $splitPoint)

# Re-touch the text of the run that follows the bookmark so it becomes a
# genuinely distinct run (placeholder swap avoids collapsing back into the
# preceding run).
$rng1b = $d.Content
$rng1b.Find.Execute("make about Kickstarter campaigns given the provided data?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng1b.Text = "X"
$rng1b.Text = "make about Kickstarter campaigns given the provided data?"

# ---------------------------------------------------------------------------
# 3. Split "world" out of the limitations paragraph into its own run and
#    change it to "realm".
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("world", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("_tempSplit", $rng2)
$rng2.Text = "realm"
$d.Bookmarks("_tempSplit").Delete()
